$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Stash the format of the F39:G39 merged cell (it has a special medium
# border that only exists on this cell) into a scratch area far outside
# the used range, BEFORE we unmerge anything. Unmerging a merged range
# that carries a "medium" left/right border causes Excel to split the
# border into two distinct new styles, which we don't want. We will
# restore the saved format after the row shift is complete.
# ------------------------------------------------------------------
$ws.Range("F39:G39").Copy()
$ws.Range("Z100:AA100").PasteSpecial(-4122)  # xlPasteFormats

# ------------------------------------------------------------------
# Unmerge the ranges that sit on / below the row where we are about to
# insert a new row, so the insert operation does not get confused by
# pre-existing merges. We purposely leave F39:G39 merged (handled via
# the stash/restore above) to avoid corrupting its border style.
# ------------------------------------------------------------------
$ws.Range("K38:N38").UnMerge()
$ws.Range("A39:E39").UnMerge()
$ws.Range("I39:N39").UnMerge()

# Insert a new blank row at position 38 (only within the used columns,
# so we do not create a full 16384-column row).
$ws.Range("A38:N38").Insert()

# Copy the formatting of row 37 (a normal product row) onto the new
# row 38 so it re-uses the exact same style indices.
$ws.Range("A37:N37").Copy()
$ws.Range("A38:N38").PasteSpecial(-4122)  # xlPasteFormats

# Re-create the merges in the same order they appear in the final file.
$ws.Range("B38:G38").Merge()
$ws.Range("H38:K38").Merge()
$ws.Range("L38:M38").Merge()
$ws.Range("K39:N39").Merge()
$ws.Range("A40:E40").Merge()
$ws.Range("F40:G40").Merge()
$ws.Range("I40:N40").Merge()

# Restore F40:G40 (originally F39:G39) formatting from the stash.
$ws.Range("Z100:AA100").Copy()
$ws.Range("F40:G40").PasteSpecial(-4122)  # xlPasteFormats

# Clean up the scratch area completely.
$ws.Range("Z100:AA100").UnMerge()
$ws.Range("Z100:AA100").Clear()

# Row heights: new product row 38 matches the other product rows (25.5),
# the footer row (now row 40) becomes 16.5.
$ws.Rows.Item(38).RowHeight = 25.5
$ws.Rows.Item(40).RowHeight = 16.5

# ------------------------------------------------------------------
# Fill in the new product row (38) values.
# ------------------------------------------------------------------
$ws.Range("A38").Value = 35
$ws.Range("B38").Value = "معطر فريدا "
$ws.Range("H38").Value = "11:0"
$ws.Range("L38").Value = -65
$ws.Range("N38").Value = "1:0"

# ------------------------------------------------------------------
# Update the totals row (now row 39): new sum of L4:L38.
# ------------------------------------------------------------------
$ws.Range("K39").Value = 1988.6400000000001

Write-Host "edit complete"
